# Append the new run-log row (row 51) to the bottom of the sheet,
# mirroring the formatting of the previous last row (row 50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 51
$prevRow = $newRow - 1

# Copy the previous row's formatting (style, borders, alignment, etc.)
# down into the new row before writing values into it.
$ws.Range("A" + $prevRow + ":H" + $prevRow).Copy()
$ws.Range("A" + $newRow + ":H" + $newRow).PasteSpecial(-4122)

$ws.Cells.Item($newRow, 1).Value = "2025-08-24 06:46:26 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-24 12:16:26 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""
